$d = $word.ActiveDocument

function Replace-Text($findText, $replText) {
    $rng = $d.Content
    $f = $rng.Find
    $f.ClearFormatting()
    $f.Text = $findText
    $ok = $f.Execute()
    if ($ok) {
        $rng.Text = $replText
    }
    return $ok
}

# 1. "the purpose of sensors is to record ... outside world, turn it ... signal, which" ->
#    "the purpose of sensors are to record ... outside world, and turn it ... signal which"
Replace-Text `
  "the purpose of sensors is to record some measurement from the outside world, turn it into an electrical signal, which is then analyzed by a device." `
  "the purpose of sensors are to record some measurement from the outside world, and turn it into an electrical signal which is then analyzed by a device."

# 2. "(ie Electrical Sensors)" -> "(i.e. Electrical Sensors)" (inside hyperlink)
Replace-Text "Uses (ie Electrical" "Uses (i.e. Electrical"

# 3. "After processing the controller can then act on the world." ->
#    "After processing, the controller can then act on the world."
Replace-Text `
  "After processing the controller can then act on the world." `
  "After processing, the controller can then act on the world."

# 4a. insert a comma right after "(Internet of things)" (the ")" sits in a plain
#     run right after the hyperlink, so use a pure insertion to avoid picking up
#     the hyperlink's character formatting).
$rng = $d.Content
$f = $rng.Find
$f.ClearFormatting()
$f.Text = ") sensors are more important"
if ($f.Execute()) {
    $insRng = $d.Range($rng.Start + 1, $rng.Start + 1)
    $insRng.InsertBefore(",")
}

# 4b. "understand how they work and use them effectively." ->
#     "understand how they work, and how to use them effectively."
Replace-Text `
  "understand how they work and use them effectively." `
  "understand how they work, and how to use them effectively."

# 5. "5 volt pin (5V) and VCC" -> "5-volt pin (5V) and VCC"
Replace-Text "5 volt pin (5V) and VCC" "5-volt pin (5V) and VCC"

# 6. "Copy and upload the code provided with this lab to the Arduino." ->
#    "Copy and upload the code provided with this lab on GitHub to the Arduino."
Replace-Text `
  "Copy and upload the code provided with this lab to the Arduino." `
  "Copy and upload the code provided with this lab on GitHub to the Arduino."

# 6b. Word stamps a "_GoBack" bookmark at the most-recent edit location;
#     recreate it right after the newly inserted "on GitHub " text.
$rng = $d.Content
$f = $rng.Find
$f.ClearFormatting()
$f.Text = "on GitHub "
if ($f.Execute()) {
    $bmRng = $d.Range($rng.End, $rng.End)
    $d.Bookmarks.Add("_GoBack", $bmRng)
}
